$wb = $excel.ActiveWorkbook

# --- Rename the trailing "Sheet3" tab to "Reward" ---
$reward = $wb.Worksheets.Item("Sheet3")
$reward.Name = "Reward"

# --- Rework the "Required" sheet's row labels for the new ride/crew model ---
$required = $wb.Worksheets.Item("Required")
$required.Range("A1").Value = "RideNum"
$required.Range("A2").Value = "CrewRequirement"

# Leave a threaded note on the new CrewRequirement row explaining the typical value
$required.Range("A2").AddCommentThreaded("This is almost always going to be 1 but sometimes you have super sorties or other such events that require more")

# Autosize the value column on the Reward sheet now that it's in active use
$reward.Columns.Item(2).AutoFit()

# --- Make "Reward" the active/selected sheet tab ---
$reward.Activate()
